# Adds the rows logged by the next day of the auto-scheduler run (2025-07-31),
# continuing on from the existing log at row 29 (A1:D29 -> A1:D58).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 30
# Post Number, Category (blank for un-categorized batch rows), Status, Scheduled Time
$data = @(
    ,(1, $null, "✅ Scheduled", "2025-07-31 14:30:00")
    ,(2, $null, "✅ Scheduled", "2025-07-31 14:48:00")
    ,(3, $null, "✅ Scheduled", "2025-07-31 15:06:00")
    ,(4, $null, "✅ Scheduled", "2025-07-31 15:24:00")
    ,(5, $null, "✅ Scheduled", "2025-07-31 15:42:00")
    ,(6, $null, "✅ Scheduled", "2025-07-31 16:00:00")
    ,(7, $null, "✅ Scheduled", "2025-07-31 16:18:00")
    ,(8, $null, "❌ Failed: The caption is too long (caused by SendMediaRequest)", "2025-07-31 16:36:00")
    ,(9, $null, "✅ Scheduled", "2025-07-31 16:54:00")
    ,(10, $null, "✅ Scheduled", "2025-07-31 17:12:00")
    ,(11, $null, "✅ Scheduled", "2025-07-31 17:30:00")
    ,(1, "Kid's Carnival", "✅ Scheduled", "2025-07-31 14:30:00")
    ,(2, "Daily Essentials", "✅ Scheduled", "2025-07-31 16:00:00")
    ,(3, "Laptops", "✅ Scheduled", "2025-07-31 17:30:00")
    ,(1, "Kid's Carnival", "✅ Scheduled", "2025-07-31 14:30:00")
    ,(2, "Daily Essentials", "✅ Scheduled", "2025-07-31 16:00:00")
    ,(3, "Laptops", "✅ Scheduled", "2025-07-31 17:30:00")
    ,(1, "Kid's Carnival", "✅ Scheduled", "2025-07-31 14:30:00")
    ,(2, "Daily Essentials", "✅ Scheduled", "2025-07-31 16:00:00")
    ,(3, "Laptops", "✅ Scheduled", "2025-07-31 17:30:00")
    ,(1, "Kid's Carnival", "✅ Scheduled", "2025-07-31 14:30:00")
    ,(2, "Daily Essentials", "✅ Scheduled", "2025-07-31 16:00:00")
    ,(3, "Laptops", "✅ Scheduled", "2025-07-31 17:30:00")
    ,(1, "Kid's Carnival", "✅ Scheduled", "2025-07-31 14:30:00")
    ,(2, "Daily Essentials", "✅ Scheduled", "2025-07-31 16:00:00")
    ,(3, "Laptops", "✅ Scheduled", "2025-07-31 17:30:00")
    ,(1, "Kid's Carnival", "✅ Scheduled", "2025-07-31 14:30:00")
    ,(2, "Daily Essentials", "✅ Scheduled", "2025-07-31 16:00:00")
    ,(3, "Laptops", "✅ Scheduled", "2025-07-31 17:30:00")
)

$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($null -ne $row[1]) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    } else {
        # Write an explicit empty string (matching the blank Category cells
        # already present in the earlier, un-categorized batch rows) rather
        # than leaving the cell completely absent. A bare "" is treated by the
        # engine as "clear the cell", so go via the quote-prefix literal and
        # then strip the resulting quote-prefix formatting back to Normal.
        $ws.Cells.Item($r, 2).Value = "'"
        $ws.Cells.Item($r, 2).Style = "Normal"
    }
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
